$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 currently holds "LOT2038 - Tecnologia de Bebidas (Indicacao de Conjunto)"
# Row 25 currently holds "LOT2028 - Tecnologia de Processos Fermentativos (Requisito fraco)"
# The commit swaps the order of these two shared strings, so the LOT2028 entry
# should now appear on row 24 and the LOT2038 entry on row 25.

$lot2038 = "LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)`n"
$lot2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"

$ws.Range("B24").Value = $lot2028
$ws.Range("C24").Value = $lot2028

$ws.Range("B25").Value = $lot2038
$ws.Range("C25").Value = $lot2038
